# Slide 2: remove the click-triggered "fade in" entrance effect that was
# applied to the title placeholder (shape id 2). This also drops the
# corresponding <p:bldP spid="2".../> build entry and lets PowerPoint
# renumber the remaining timing node ids, matching the author's
# "minor changes on slide2" edit.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$mainSeq = $s.TimeLine.MainSequence

for ($i = $mainSeq.Count; $i -ge 1; $i--) {
    $effect = $mainSeq.Item($i)
    if ($effect.Shape.Id -eq 2) {
        $effect.Delete()
    }
}
